$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.5714285714285714; C = 0.6153846153846154; D = 0.5925925925925927; E = 13 }
    3  = @{ C = 0.4545454545454545; D = 0.4761904761904762; E = 11 }
    4  = @{ B = 0.5416666666666666; C = 0.5416666666666666; D = 0.5416666666666666; E = 0.5416666666666666 }
    5  = @{ B = 0.5357142857142857; C = 0.534965034965035;  D = 0.5343915343915344 }
    6  = @{ B = 0.5386904761904762; C = 0.5416666666666666; D = 0.539241622574956 }
    7  = @{ B = 0.5454545454545454; C = 0.4615384615384616; D = 0.4999999999999999; E = 13 }
    8  = @{ B = 0.4615384615384616; C = 0.5454545454545454; D = 0.4999999999999999; E = 11 }
    9  = @{ B = 0.5; C = 0.5; D = 0.5; E = 0.5 }
    10 = @{ B = 0.5034965034965035; C = 0.5034965034965035; D = 0.4999999999999999 }
    11 = @{ B = 0.506993006993007;  C = 0.5;                D = 0.4999999999999998 }
    12 = @{ B = 0.375; C = 0.2307692307692308; D = 0.2857142857142857; E = 13 }
    13 = @{ B = 0.375; C = 0.5454545454545454; D = 0.4444444444444444; E = 11 }
    14 = @{ B = 0.375; C = 0.375; D = 0.375; E = 0.375 }
    15 = @{ B = 0.375; C = 0.3881118881118881; D = 0.3650793650793651 }
    16 = @{ B = 0.375; C = 0.375;              D = 0.3584656084656084 }
    17 = @{ B = 0.5333333333333333; C = 0.6153846153846154; D = 0.5714285714285715; E = 13 }
    18 = @{ B = 0.4444444444444444; C = 0.3636363636363636; D = 0.4; E = 11 }
    19 = @{ B = 0.5; C = 0.5; D = 0.5; E = 0.5 }
    20 = @{ B = 0.4888888888888889; C = 0.4895104895104895; D = 0.4857142857142858 }
    21 = @{ B = 0.4925925925925926; C = 0.5;                D = 0.4928571428571429 }
    22 = @{ B = 0.6; C = 0.4615384615384616; D = 0.5217391304347826; E = 13 }
    23 = @{ B = 0.5; C = 0.6363636363636364; D = 0.5600000000000001; E = 11 }
    24 = @{ B = 0.5416666666666666; C = 0.5416666666666666; D = 0.5416666666666666; E = 0.5416666666666666 }
    25 = @{ B = 0.55; C = 0.548951048951049; D = 0.5408695652173914 }
    26 = @{ B = 0.5541666666666667; C = 0.5416666666666666; D = 0.5392753623188405 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
